# Update countries & provincias Spain
# Applies updated COVID-19 case figures and the resulting re-ranking of a few
# countries whose totals crossed each other, plus the refreshed timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 17:22"

# Row 4: updated figures
$ws.Range("B4").Value = 475322
$ws.Range("C4").Value = 6756
$ws.Range("E4").Value = 432217

# Row 19: updated figures
$ws.Range("B19").Value = 13520
$ws.Range("C19").Value = 276
$ws.Range("E19").Value = 7137

# Row 30: Dinamarca -> Polonia
$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 5955
$ws.Range("C30").Value = 380
$ws.Range("D30").Value = 318
$ws.Range("E30").Value = 5456
$ws.Range("F30").Value = 160
$ws.Range("G30").Value = 7
$ws.Range("H30").Value = 181

# Row 31: Polonia -> Dinamarca
$ws.Range("A31").Value = "Dinamarca"
$ws.Range("B31").Value = 5819
$ws.Range("C31").Value = 184
$ws.Range("D31").Value = 1773
$ws.Range("E31").Value = 3799
$ws.Range("F31").Value = 113
$ws.Range("G31").Value = 10
$ws.Range("H31").Value = 247

# Row 53: Sudafrica -> Grecia
$ws.Range("A53").Value = "Grecia"
$ws.Range("B53").Value = 2011
$ws.Range("C53").Value = 56
$ws.Range("D53").Value = 269
$ws.Range("E53").Value = 1652
$ws.Range("F53").Value = 79
$ws.Range("G53").Value = 3
$ws.Range("H53").Value = 90

# Row 54: Bielorrusia -> Sudafrica
$ws.Range("A54").Value = "Sudafrica"
$ws.Range("B54").Value = 2003
$ws.Range("C54").Value = 69
$ws.Range("D54").Value = 410
$ws.Range("E54").Value = 1569
$ws.Range("F54").Value = 7
$ws.Range("G54").Value = 6
$ws.Range("H54").Value = 24

# Row 55: Grecia -> Bielorrusia
$ws.Range("A55").Value = "Bielorrusia"
$ws.Range("B55").Value = 1981
$ws.Range("C55").Value = 495
$ws.Range("D55").Value = 169
$ws.Range("E55").Value = 1793
$ws.Range("F55").Value = 72
$ws.Range("G55").Value = 3
$ws.Range("H55").Value = 19

# Row 56: updated figures
$ws.Range("E56").Value = 1448
$ws.Range("G56").Value = 2
$ws.Range("H56").Value = 81

# Row 70: Hong Kong -> Azerbaiyan
$ws.Range("A70").Value = "Azerbaiyan"
$ws.Range("B70").Value = 991
$ws.Range("C70").Value = 65
$ws.Range("D70").Value = 159
$ws.Range("E70").Value = 822
$ws.Range("F70").Value = 27
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 10

# Row 71: Armenia -> Hong Kong
$ws.Range("A71").Value = "Hong Kong"
$ws.Range("B71").Value = 990
$ws.Range("D71").Value = 309
$ws.Range("E71").Value = 677
$ws.Range("F71").Value = 15
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 4

# Row 72: Azerbaiyan -> Armenia
$ws.Range("A72").Value = "Armenia"
$ws.Range("B72").Value = 937
$ws.Range("C72").Value = 16
$ws.Range("D72").Value = 149
$ws.Range("E72").Value = 776
$ws.Range("F72").Value = 30
$ws.Range("G72").Value = 2
$ws.Range("H72").Value = 12

# Row 87: Costa Rica -> Cuba
$ws.Range("A87").Value = "Cuba"
$ws.Range("B87").Value = 564
$ws.Range("C87").Value = 49
$ws.Range("D87").Value = 51
$ws.Range("E87").Value = 498
$ws.Range("F87").Value = 15
$ws.Range("H87").Value = 15

# Row 88: Afganistan -> Costa Rica
$ws.Range("A88").Value = "Costa Rica"
$ws.Range("B88").Value = 539
$ws.Range("C88").Value = 0
$ws.Range("D88").Value = 30
$ws.Range("E88").Value = 506
$ws.Range("F88").Value = 13
$ws.Range("H88").Value = 3

# Row 89: Cuba -> Afganistan
$ws.Range("A89").Value = "Afganistan"
$ws.Range("B89").Value = 521
$ws.Range("C89").Value = 37
$ws.Range("D89").Value = 32
$ws.Range("E89").Value = 474
$ws.Range("F89").Value = 0

# Row 114: Guinea -> Isla de Man
$ws.Range("A114").Value = "Isla de Man"
$ws.Range("B114").Value = 201
$ws.Range("C114").Value = 11
$ws.Range("D114").Value = 103
$ws.Range("E114").Value = 97
$ws.Range("F114").Value = 11
$ws.Range("H114").Value = 1

# Row 115: Mayotte -> Guinea
$ws.Range("A115").Value = "Guinea"
$ws.Range("B115").Value = 194
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 11
$ws.Range("E115").Value = 183
$ws.Range("F115").Value = 0
$ws.Range("H115").Value = 0

# Row 116: Sri Lanka -> Mayotte
$ws.Range("A116").Value = "Mayotte"
$ws.Range("B116").Value = 191
$ws.Range("C116").Value = 7
$ws.Range("D116").Value = 50
$ws.Range("E116").Value = 139
$ws.Range("F116").Value = 4
$ws.Range("H116").Value = 2

# Row 117: Isla de Man -> Sri Lanka
$ws.Range("A117").Value = "Sri Lanka"
$ws.Range("D117").Value = 54
$ws.Range("E117").Value = 129
$ws.Range("F117").Value = 5
$ws.Range("H117").Value = 7

# Row 153: Guyana -> Liberia
$ws.Range("A153").Value = "Liberia"
$ws.Range("C153").Value = 6
$ws.Range("D153").Value = 3
$ws.Range("E153").Value = 29
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 5

# Row 154: Guinea-Bisau -> Guyana
$ws.Range("A154").Value = "Guyana"
$ws.Range("B154").Value = 37
$ws.Range("D154").Value = 8
$ws.Range("E154").Value = 23
$ws.Range("F154").Value = 4
$ws.Range("H154").Value = 6

# Row 155: Eritrea -> Guinea-Bisau
$ws.Range("A155").Value = "Guinea-Bisau"
$ws.Range("B155").Value = 36
$ws.Range("E155").Value = 36

# Row 156: Guam -> Benin
$ws.Range("A156").Value = "Benin"
$ws.Range("B156").Value = 35
$ws.Range("C156").Value = 9
$ws.Range("D156").Value = 5
$ws.Range("E156").Value = 29

# Row 157: Tanzania -> Eritrea
$ws.Range("A157").Value = "Eritrea"
$ws.Range("B157").Value = 33
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 0
$ws.Range("E157").Value = 33
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 0

# Row 158: San Martin (Parte Francesa) -> Guam
$ws.Range("A158").Value = "Guam"
$ws.Range("D158").Value = 0
$ws.Range("E158").Value = 31
$ws.Range("F158").Value = 0
$ws.Range("H158").Value = 1

# Row 159: Liberia -> Tanzania
$ws.Range("A159").Value = "Tanzania"
$ws.Range("B159").Value = 32
$ws.Range("C159").Value = 7
$ws.Range("D159").Value = 5
$ws.Range("G159").Value = 2
$ws.Range("H159").Value = 3

# Row 160: Haiti -> San Martin (Parte Francesa)
$ws.Range("A160").Value = "San Martin (Parte Francesa)"
$ws.Range("B160").Value = 32
$ws.Range("D160").Value = 11
$ws.Range("E160").Value = 19
$ws.Range("F160").Value = 5

# Row 161: Birmania -> Haiti
$ws.Range("A161").Value = "Haiti"
$ws.Range("B161").Value = 30
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 0
$ws.Range("E161").Value = 28
$ws.Range("H161").Value = 2

# Row 162: Benin -> Birmania
$ws.Range("A162").Value = "Birmania"
$ws.Range("B162").Value = 27
$ws.Range("C162").Value = 4
$ws.Range("D162").Value = 2
$ws.Range("E162").Value = 22
$ws.Range("H162").Value = 3

# Row 172: Fiyi -> Laos
$ws.Range("A172").Value = "Laos"
$ws.Range("C172").Value = 0

# Row 173: Laos -> Fiyi
$ws.Range("A173").Value = "Fiyi"
$ws.Range("C173").Value = 1

# Row 185: San Cristobal y Nieves -> Seychelles
$ws.Range("A185").Value = "Seychelles"

# Row 186: Seychelles -> San Cristobal y Nieves
$ws.Range("A186").Value = "San Cristobal y Nieves"
